$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1,1)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "87÷2=43, 1"

$c = $t.Cell(1,2)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "33÷9=3, 6"

$c = $t.Cell(1,3)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "51÷5=10, 1"

$c = $t.Cell(1,4)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "47÷2=23, 1"

$c = $t.Cell(1,5)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "97÷2=48, 1"

$c = $t.Cell(5,1)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "71÷2=35, 1"

$c = $t.Cell(5,2)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "71÷5=14, 1"

$c = $t.Cell(5,3)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "55÷6=9, 1"

$c = $t.Cell(5,4)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "99÷5=19, 4"

$c = $t.Cell(5,5)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "40÷8=5, 0"

$c = $t.Cell(9,1)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "58÷2=29, 0"

$c = $t.Cell(9,2)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "13÷8=1, 5"

$c = $t.Cell(9,3)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "53÷6=8, 5"

$c = $t.Cell(9,4)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "35÷2=17, 1"

$c = $t.Cell(9,5)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "79÷7=11, 2"

$c = $t.Cell(13,1)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "36÷5=7, 1"

$c = $t.Cell(13,2)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "30÷9=3, 3"

$c = $t.Cell(13,3)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "51÷9=5, 6"

$c = $t.Cell(13,4)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "28÷3=9, 1"

$c = $t.Cell(13,5)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "49÷3=16, 1"

$c = $t.Cell(17,1)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "90÷9=10, 0"

$c = $t.Cell(17,2)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "60÷4=15, 0"

$c = $t.Cell(17,3)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "20÷3=6, 2"

$c = $t.Cell(17,4)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "35÷7=5, 0"

$c = $t.Cell(17,5)
$r = $c.Range
$d.Range($r.Start, $r.End - 1).Text = "55÷3=18, 1"

